$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.5586792620790276
$ws.Range("J2").Value = 0.5586792620790275
$ws.Range("O2").Value = 0.7091726973716084
$ws.Range("P2").Value = 0.7091726973716084
$ws.Range("S2").Value = 0.3962000792541637
$ws.Range("T2").Value = 0.3962000792541637

# Row 3
$ws.Range("I3").Value = 0.5586792620790276
$ws.Range("J3").Value = 0.5586792620790275
$ws.Range("M3").Value = 0.7003526666666667
$ws.Range("N3").Value = 2.101058
$ws.Range("O3").Value = 0.2908273026283917
$ws.Range("P3").Value = 0.2908273026283917
$ws.Range("Q3").Value = 0.08644266203955557
$ws.Range("R3").Value = 0.777983958356
$ws.Range("S3").Value = 0.1624791828248639
$ws.Range("T3").Value = 0.1624791828248639

# Row 4
$ws.Range("G4").Value = 0.09749966666666667
$ws.Range("H4").Value = 0.292499
$ws.Range("I4").Value = 0.4413207379209724
$ws.Range("J4").Value = 0.4413207379209724
$ws.Range("O4").Value = 0.7091726973716084
$ws.Range("P4").Value = 0.7091726973716084
$ws.Range("Q4").Value = 0.1665086307377778
$ws.Range("R4").Value = 1.49857767664
$ws.Range("S4").Value = 0.3129726181174446
$ws.Range("T4").Value = 0.3129726181174446

# Row 5
$ws.Range("G5").Value = 0.09749966666666667
$ws.Range("H5").Value = 0.292499
$ws.Range("I5").Value = 0.4413207379209724
$ws.Range("J5").Value = 0.4413207379209724
$ws.Range("M5").Value = 0.7003526666666667
$ws.Range("N5").Value = 2.101058
$ws.Range("O5").Value = 0.2908273026283917
$ws.Range("P5").Value = 0.2908273026283917
$ws.Range("Q5").Value = 0.06828415154911112
$ws.Range("R5").Value = 0.614557363942
$ws.Range("S5").Value = 0.1283481198035278
$ws.Range("T5").Value = 0.1283481198035278
